$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '68.661.96'
Set-TextValue $ws.Range('D3') '3.860.07'
Set-TextValue $ws.Range('E3') '  -1.40%  '
Set-TextValue $ws.Range('D4') '1.00'
Set-TextValue $ws.Range('E4') '  -0.05%  '
Set-TextValue $ws.Range('D5') '602.75'
Set-TextValue $ws.Range('E5') '  -0.52%  '
Set-TextValue $ws.Range('D6') '168.74'
Set-TextValue $ws.Range('E6') '  +2.69%  '
Set-TextValue $ws.Range('D7') '3.860.67'
Set-TextValue $ws.Range('E7') '  -1.32%  '
Set-TextValue $ws.Range('E8') '  -0.02%  '
Set-TextValue $ws.Range('E9') '  -0.26%  '
Set-TextValue $ws.Range('E10') '  -0.76%  '
Set-TextValue $ws.Range('D11') '6.36'
Set-TextValue $ws.Range('E11') '  -0.85%  '
Set-TextValue $ws.Range('E12') '  +0.32%  '
Set-TextValue $ws.Range('D13') '0.0000251'
Set-TextValue $ws.Range('E13') '  +2.23%  '
Set-TextValue $ws.Range('D14') '37.60'
Set-TextValue $ws.Range('E14') '  +1.22%  '
Set-TextValue $ws.Range('D15') '4.512.18'
Set-TextValue $ws.Range('E15') '  -1.23%  '
Set-TextValue $ws.Range('D16') '3.872.77'
Set-TextValue $ws.Range('E16') '  -0.52%  '
Set-TextValue $ws.Range('D17') '68.798.06'
Set-TextValue $ws.Range('E17') '  -0.52%  '
Set-TextValue $ws.Range('D18') '7.58'
Set-TextValue $ws.Range('E18') '  +1.23%  '
Set-TextValue $ws.Range('D19') '18.31'
Set-TextValue $ws.Range('E19') '  +6.70%  '
Set-TextValue $ws.Range('E20') '  -1.18%  '
Set-TextValue $ws.Range('D21') '10.88'
Set-TextValue $ws.Range('E21') '  -2.99%  '
Set-TextValue $ws.Range('D22') '475.28'
Set-TextValue $ws.Range('E22') '  -2.68%  '
Set-TextValue $ws.Range('D23') '0.740'
Set-TextValue $ws.Range('E23') '  +2.17%  '
Set-TextValue $ws.Range('D24') '0.0000161'
Set-TextValue $ws.Range('E24') '  -2.59%  '
Set-TextValue $ws.Range('D25') '84.83'
Set-TextValue $ws.Range('E25') '  +0.47%  '
Set-TextValue $ws.Range('E26') '  -0.15%  '
Set-TextValue $ws.Range('D27') '12.45'
Set-TextValue $ws.Range('E27') '  +2.36%  '
Set-TextValue $ws.Range('D28') '10.12'
Set-TextValue $ws.Range('E28') '  +0.18%  '
Set-TextValue $ws.Range('E29') '  -0.01%  '
Set-TextValue $ws.Range('E30') '  +0.66%  '
Set-TextValue $ws.Range('D31') '4.011.51'
Set-TextValue $ws.Range('E31') '  -1.36%  '
Set-TextValue $ws.Range('D32') '7.81'
Set-TextValue $ws.Range('E32') '  -1.11%  '
Set-TextValue $ws.Range('E33') '  -2.28%  '
Set-TextValue $ws.Range('D34') '31.22'
Set-TextValue $ws.Range('E34') '  -3.55%  '
Set-TextValue $ws.Range('D35') '3.829.47'
Set-TextValue $ws.Range('E35') '  -0.74%  '
Set-TextValue $ws.Range('E36') '  -1.01%  '
Set-TextValue $ws.Range('B37') 'Kaspa'
Set-TextValue $ws.Range('C37') 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D37') '0.142'
Set-TextValue $ws.Range('E37') '  +1.21%  '
Set-TextValue $ws.Range('D38') '6.02'
Set-TextValue $ws.Range('E38') '  +1.16%  '
Set-TextValue $ws.Range('B39') 'dogwifhat'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D39') '3.35'
Set-TextValue $ws.Range('E39') '  +10.29%  '
Set-TextValue $ws.Range('D40') '1.01'
Set-TextValue $ws.Range('E40') '  -2.71%  '
Set-TextValue $ws.Range('D41') '1.00'
Set-TextValue $ws.Range('E41') '  +0.06%  '
Set-TextValue $ws.Range('D42') '0.318'
Set-TextValue $ws.Range('E42') '  -0.91%  '
Set-TextValue $ws.Range('D43') '2.02'
Set-TextValue $ws.Range('E43') '  +0.78%  '
Set-TextValue $ws.Range('D44') '431.01'
Set-TextValue $ws.Range('E44') '  -1.84%  '
Set-TextValue $ws.Range('D45') '47.66'
Set-TextValue $ws.Range('E45') '  -1.63%  '
Set-TextValue $ws.Range('E46') '  +0.00%  '
Set-TextValue $ws.Range('E47') '  +2.14%  '
Set-TextValue $ws.Range('E48') '  +14.22%  '
Set-TextValue $ws.Range('E49') '  +1.22%  '
Set-TextValue $ws.Range('D50') '142.25'
Set-TextValue $ws.Range('E50') '  +0.62%  '
Set-TextValue $ws.Range('D51') '39.34'
Set-TextValue $ws.Range('E51') '  +0.74%  '
